$wb = $excel.ActiveWorkbook

# --- Sheet: Resumen ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z1"
$wsResumen.Range("C2").Value = 504.6388420369769

# --- Sheet: Solucion ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsSolucion.Range("B4").Value = "S011"
$wsSolucion.Range("B5").Value = "S002"
$wsSolucion.Range("B7").Value = "S031"
$wsSolucion.Range("B13").Value = "S033"
$wsSolucion.Range("B14").Value = "S024"
$wsSolucion.Range("B15").Value = "S004"
$wsSolucion.Range("B22").Value = "S006"
$wsSolucion.Range("B23").Value = "S026"
$wsSolucion.Range("B24").Value = "S016"
$wsSolucion.Range("B25").Value = "S036"
$wsSolucion.Range("B26").Value = "S027"
$wsSolucion.Range("B28").Value = "S028"
$wsSolucion.Range("B30").Value = "S037"
$wsSolucion.Range("B31").Value = "S018"
$wsSolucion.Range("B33").Value = "S008"
$wsSolucion.Range("B36").Value = "S039"
$wsSolucion.Range("B37").Value = "S019"
$wsSolucion.Range("B40").Value = "S040"
$wsSolucion.Range("B41").Value = "S020"

# --- Sheet: Metricas ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 504.6388420369769
$wsMetricas.Range("B3").Value = 504.5936641799115
